$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = 1.27
$ws.Range("E2").Value = 1.33

$ws.Range("B4").Value = 1.45
$ws.Range("C4").Value = 1.4
$ws.Range("F4").Value = 1.08

$ws.Range("B5").Value = 1.42
$ws.Range("D5").Value = 1.31
$ws.Range("F5").Value = 1.05

$ws.Range("D6").Value = 1.52
